# Update the Gnai2-Adra2a LR-pairs sheet with refreshed TPM-derived
# ligand/receptor expression stats (per commit "update scripts wuth new tpm").
# Only the numeric expression/specificity/edge-weight columns (G,H,I,J,M,N,O,P,Q,R,S,T)
# change; identifiers in columns A-F, K, L are untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 82.98768099999999
$ws.Range("H2").Value = 248.963043
$ws.Range("I2").Value = 0.4489504115427952
$ws.Range("J2").Value = 0.4489504115427952
$ws.Range("M2").Value = 0.1999913333333333
$ws.Range("N2").Value = 0.599974
$ws.Range("O2").Value = 0.03856011821502105
$ws.Range("P2").Value = 0.03856011821502105
$ws.Range("Q2").Value = 16.59681697343133
$ws.Range("R2").Value = 149.371352760882
$ws.Range("S2").Value = 0.01731158094177254
$ws.Range("T2").Value = 0.01731158094177254
$ws.Range("G3").Value = 82.98768099999999
$ws.Range("H3").Value = 248.963043
$ws.Range("I3").Value = 0.4489504115427952
$ws.Range("J3").Value = 0.4489504115427952
$ws.Range("O3").Value = 0.164275600079283
$ws.Range("P3").Value = 0.164275600079283
$ws.Range("Q3").Value = 70.70652772673232
$ws.Range("R3").Value = 636.358749540591
$ws.Range("S3").Value = 0.07375159826203377
$ws.Range("T3").Value = 0.07375159826203377
$ws.Range("G4").Value = 82.98768099999999
$ws.Range("H4").Value = 248.963043
$ws.Range("I4").Value = 0.4489504115427952
$ws.Range("J4").Value = 0.4489504115427952
$ws.Range("O4").Value = 0.7971642817056959
$ws.Range("P4").Value = 0.797164281705696
$ws.Range("Q4").Value = 343.1107137029576
$ws.Range("R4").Value = 3087.996423326619
$ws.Range("S4").Value = 0.3578872323389889
$ws.Range("T4").Value = 0.3578872323389889
$ws.Range("G5").Value = 63.14058933333333
$ws.Range("I5").Value = 0.3415807409566563
$ws.Range("J5").Value = 0.3415807409566563
$ws.Range("M5").Value = 0.1999913333333333
$ws.Range("N5").Value = 0.599974
$ws.Range("O5").Value = 0.03856011821502105
$ws.Range("P5").Value = 0.03856011821502105
$ws.Range("Q5").Value = 12.62757064822578
$ws.Range("R5").Value = 113.648135834032
$ws.Range("S5").Value = 0.01317139375126315
$ws.Range("T5").Value = 0.01317139375126315
$ws.Range("G6").Value = 63.14058933333333
$ws.Range("I6").Value = 0.3415807409566563
$ws.Range("J6").Value = 0.3415807409566563
$ws.Range("O6").Value = 0.164275600079283
$ws.Range("P6").Value = 0.164275600079283
$ws.Range("Q6").Value = 53.7965608459351
$ws.Range("S6").Value = 0.05611338119618085
$ws.Range("T6").Value = 0.05611338119618085
$ws.Range("G7").Value = 63.14058933333333
$ws.Range("I7").Value = 0.3415807409566563
$ws.Range("J7").Value = 0.3415807409566563
$ws.Range("O7").Value = 0.7971642817056959
$ws.Range("P7").Value = 0.797164281705696
$ws.Range("S7").Value = 0.2722959660092123
$ws.Range("T7").Value = 0.2722959660092124
$ws.Range("I8").Value = 0.2094688475005485
$ws.Range("J8").Value = 0.2094688475005485
$ws.Range("M8").Value = 0.1999913333333333
$ws.Range("N8").Value = 0.599974
$ws.Range("O8").Value = 0.03856011821502105
$ws.Range("P8").Value = 0.03856011821502105
$ws.Range("Q8").Value = 7.743652827169334
$ws.Range("R8").Value = 69.69287544452401
$ws.Range("S8").Value = 0.008077143521985367
$ws.Range("T8").Value = 0.008077143521985367
$ws.Range("I9").Value = 0.2094688475005485
$ws.Range("J9").Value = 0.2094688475005485
$ws.Range("O9").Value = 0.164275600079283
$ws.Range("P9").Value = 0.164275600079283
$ws.Range("S9").Value = 0.03441062062106843
$ws.Range("T9").Value = 0.03441062062106843
$ws.Range("I10").Value = 0.2094688475005485
$ws.Range("J10").Value = 0.2094688475005485
$ws.Range("O10").Value = 0.7971642817056959
$ws.Range("P10").Value = 0.797164281705696
$ws.Range("Q10").Value = 160.0867354536287
$ws.Range("S10").Value = 0.1669810833574947
$ws.Range("T10").Value = 0.1669810833574947

Write-Output "Applied 81 cell updates to worksheet."
